$wb = $excel.ActiveWorkbook

$wsCommon = $wb.Worksheets.Item("common")
$wsSimulator = $wb.Worksheets.Item("simulator")

# Row heights: rows 2-7 on "common" grow from 14.25 to 27 (wrapped/taller rows)
$wsCommon.Range("A2:D7").RowHeight = 27

# New localisation rows for the ACTIVE / CLOSED status strings
$wsCommon.Cells.Item(8, 1).Value = "status.active"
$wsCommon.Cells.Item(8, 2).Value = "开启"
$wsCommon.Cells.Item(8, 3).Value = "ACTIVE"
$wsCommon.Cells.Item(8, 4).Value = "ACTIVE"

$wsCommon.Cells.Item(9, 1).Value = "status.closed"
$wsCommon.Cells.Item(9, 2).Value = "关闭"
$wsCommon.Cells.Item(9, 3).Value = "CLOSED"
$wsCommon.Cells.Item(9, 4).Value = "CLOSED"

# Match the formatting used by column D on the existing rows (style index 2)
$wsCommon.Range("D7").Copy()
$wsCommon.Range("D8").PasteSpecial(-4122)
$wsCommon.Range("D9").PasteSpecial(-4122)

# New rows also use the taller row height
$wsCommon.Range("A8:D9").RowHeight = 27

# Theme: lighten the window background colour back to white
$theme = $wb.Theme
$colorScheme = $theme.ThemeColorScheme
$colorScheme.Colors(2).RGB = 16777215

# Restore cursor/selection positions on both sheets
$wsCommon.Activate()
$wsCommon.Range("C15").Select()

$wsSimulator.Activate()
$wsSimulator.Range("B13").Select()
